# Generate Report for Handoff
# Update status from "In Translation" to "Ready for handoff" and refresh
# the handoff timestamps across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns + overall latest handoff date
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-06-17 10:06:11"

# zh-cn detail sheet: Status + Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-17 10:06:08"

# de-de detail sheet: Status + Latest Handoff Datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-17 10:06:11"
